$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -4
$ws.Range("J25").Value = -4.3
$ws.Range("K25").Value = -1
$ws.Range("I26").Value = -4.9
$ws.Range("J26").Value = -1.6
$ws.Range("K26").Value = -0.4
$ws.Range("H27").Value = -3.9
$ws.Range("I27").Value = -0.6
$ws.Range("J27").Value = 0.6
$ws.Range("K27").Value = 0.6
$ws.Range("G28").Value = -3.7
$ws.Range("H28").Value = -0.4
$ws.Range("I28").Value = 0.8
$ws.Range("J28").Value = 0.8
$ws.Range("K28").Value = 0.4
$ws.Range("F29").Value = -4.4
$ws.Range("G29").Value = -1.1
$ws.Range("H29").Value = 0.09999999999999998
$ws.Range("I29").Value = 0.09999999999999998
$ws.Range("J29").Value = -0.3
$ws.Range("K29").Value = -0.7
$ws.Range("E30").Value = -3.4
$ws.Range("F30").Value = -0.09999999999999998
$ws.Range("G30").Value = 1.1
$ws.Range("H30").Value = 1.1
$ws.Range("I30").Value = 0.7
$ws.Range("J30").Value = 0.3
$ws.Range("K30").Value = 1.9
$ws.Range("D31").Value = -4.5
$ws.Range("E31").Value = -1.2
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = -0.4
$ws.Range("I31").Value = -0.8
$ws.Range("J31").Value = 0.7999999999999999
$ws.Range("K31").Value = 0
$ws.Range("C32").Value = -4.7
$ws.Range("D32").Value = -1.4
$ws.Range("E32").Value = -0.2000000000000001
$ws.Range("F32").Value = -0.2000000000000001
$ws.Range("G32").Value = -0.6000000000000001
$ws.Range("H32").Value = -1
$ws.Range("I32").Value = 0.5999999999999999
$ws.Range("J32").Value = -0.2000000000000001
$ws.Range("K32").Value = -0.7000000000000001
$ws.Range("B33").Value = -6.699999999999999
$ws.Range("C33").Value = -3.4
$ws.Range("D33").Value = -2.2
$ws.Range("E33").Value = -2.2
$ws.Range("F33").Value = -2.6
$ws.Range("G33").Value = -3
$ws.Range("H33").Value = -1.4
$ws.Range("I33").Value = -2.2
$ws.Range("J33").Value = -2.7
$ws.Range("K33").Value = -2.9
$ws.Range("B34").Value = 3.3
$ws.Range("C34").Value = 4.5
$ws.Range("D34").Value = 4.5
$ws.Range("E34").Value = 4.1
$ws.Range("F34").Value = 3.7
$ws.Range("G34").Value = 5.3
$ws.Range("H34").Value = 4.5
$ws.Range("I34").Value = 4
$ws.Range("J34").Value = 3.8
$ws.Range("K34").Value = 3.2
$ws.Range("B35").Value = 1.2
$ws.Range("C35").Value = 1.2
$ws.Range("D35").Value = 0.8
$ws.Range("E35").Value = 0.4
$ws.Range("F35").Value = 2
$ws.Range("G35").Value = 1.2
$ws.Range("H35").Value = 0.7
$ws.Range("I35").Value = 0.5
$ws.Range("J35").Value = -0.09999999999999998
$ws.Range("K35").Value = 1.2
$ws.Range("C36").Value = -0.4
$ws.Range("D36").Value = -0.8
$ws.Range("E36").Value = 0.7999999999999999
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = -0.5
$ws.Range("H36").Value = -0.7
$ws.Range("I36").Value = -1.3
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = -0.8
$ws.Range("B37").Value = -0.4
$ws.Range("C37").Value = -0.8
$ws.Range("D37").Value = 0.7999999999999999
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = -0.5
$ws.Range("G37").Value = -0.7
$ws.Range("H37").Value = -1.3
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = -0.8
$ws.Range("K37").Value = -0.6
$ws.Range("B38").Value = -0.4
$ws.Range("C38").Value = 1.2
$ws.Range("D38").Value = 0.4
$ws.Range("E38").Value = -0.1
$ws.Range("F38").Value = -0.3
$ws.Range("G38").Value = -0.8999999999999999
$ws.Range("H38").Value = 0.4
$ws.Range("I38").Value = -0.4
$ws.Range("J38").Value = -0.2
$ws.Range("K38").Value = -0.5
$ws.Range("B39").Value = 1.6
$ws.Range("C39").Value = 0.8
$ws.Range("D39").Value = 0.3
$ws.Range("E39").Value = 0.1
$ws.Range("F39").Value = -0.4999999999999999
$ws.Range("G39").Value = 0.8
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0.2
$ws.Range("J39").Value = -0.09999999999999998
$ws.Range("K39").Value = 0.3
$ws.Range("B40").Value = -0.7999999999999999
$ws.Range("C40").Value = -1.3
$ws.Range("D40").Value = -1.5
$ws.Range("E40").Value = -2.1
$ws.Range("F40").Value = -0.7999999999999999
$ws.Range("G40").Value = -1.6
$ws.Range("H40").Value = -1.4
$ws.Range("I40").Value = -1.7
$ws.Range("J40").Value = -1.3
$ws.Range("K40").Value = -0.4999999999999999
$ws.Range("B41").Value = -0.5
$ws.Range("C41").Value = -0.7
$ws.Range("D41").Value = -1.3
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = -0.8
$ws.Range("G41").Value = -0.6
$ws.Range("H41").Value = -0.8999999999999999
$ws.Range("I41").Value = -0.5
$ws.Range("J41").Value = 0.3
$ws.Range("K41").Value = 0.6
$ws.Range("B42").Value = -0.2
$ws.Range("C42").Value = -0.7999999999999999
$ws.Range("D42").Value = 0.5
$ws.Range("E42").Value = -0.3
$ws.Range("F42").Value = -0.1
$ws.Range("G42").Value = -0.4
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0.8
$ws.Range("J42").Value = 1.1
$ws.Range("K42").Value = -1
$ws.Range("B43").Value = -0.6
$ws.Range("C43").Value = 0.7
$ws.Range("D43").Value = -0.1
$ws.Range("E43").Value = 0.1
$ws.Range("F43").Value = -0.2
$ws.Range("G43").Value = 0.2
$ws.Range("H43").Value = 1
$ws.Range("I43").Value = 1.3
$ws.Range("J43").Value = -0.8
$ws.Range("K43").Value = 0.7999999999999999
$ws.Range("B44").Value = 1.3
$ws.Range("C44").Value = 0.4999999999999999
$ws.Range("D44").Value = 0.7
$ws.Range("E44").Value = 0.4
$ws.Range("F44").Value = 0.7999999999999999
$ws.Range("G44").Value = 1.6
$ws.Range("H44").Value = 1.9
$ws.Range("I44").Value = -0.2000000000000001
$ws.Range("J44").Value = 1.4
$ws.Range("B45").Value = -0.8
$ws.Range("C45").Value = -0.6
$ws.Range("D45").Value = -0.8999999999999999
$ws.Range("E45").Value = -0.5
$ws.Range("F45").Value = 0.3
$ws.Range("G45").Value = 0.6
$ws.Range("H45").Value = -1.5
$ws.Range("I45").Value = 0.09999999999999998
$ws.Range("B46").Value = 0.2
$ws.Range("C46").Value = -0.09999999999999998
$ws.Range("D46").Value = 0.3
$ws.Range("E46").Value = 1.1
$ws.Range("F46").Value = 1.4
$ws.Range("G46").Value = -0.7
$ws.Range("H46").Value = 0.8999999999999999
$ws.Range("B47").Value = -0.3
$ws.Range("C47").Value = 0.1
$ws.Range("D47").Value = 0.9
$ws.Range("E47").Value = 1.2
$ws.Range("F47").Value = -0.9
$ws.Range("G47").Value = 0.7
$ws.Range("B48").Value = 0.4
$ws.Range("C48").Value = 1.2
$ws.Range("D48").Value = 1.5
$ws.Range("E48").Value = -0.6000000000000001
$ws.Range("F48").Value = 1
$ws.Range("B49").Value = 0.8
$ws.Range("C49").Value = 1.1
$ws.Range("D49").Value = -1
$ws.Range("E49").Value = 0.6
$ws.Range("B50").Value = 0.2999999999999999
$ws.Range("C50").Value = -1.8
$ws.Range("D50").Value = -0.2000000000000001
$ws.Range("B51").Value = -2.1
$ws.Range("C51").Value = -0.5
$ws.Range("B52").Value = 1.6
